# Update the "想去人数" (interested count) figures for a few rows on both
# the "展览" and "全部类型" worksheets, reflecting newly generated data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 1662
    $ws.Range("F6").Value = 438
    $ws.Range("F9").Value = 558
}
